# Inserts a new weekly record at row 91 of the "Haba" (broad bean) sheet,
# shifting the existing rows 91:172 down to 92:173, then populates the
# new row with the new price-report values for that week.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 91:172 down by one row (creates an empty row 91).
$ws.Rows("91:91").Insert()

# Populate the newly inserted row 91 with the new record.
# Columns that stay constant across this market/category block:
$ws.Cells.Item(91, 1).Value = 9
$ws.Cells.Item(91, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(91, 3).Value = "Metropolitana"
$ws.Cells.Item(91, 4).Value = 44512
$ws.Cells.Item(91, 5).Value = 13
$ws.Cells.Item(91, 6).Value = 100112026
$ws.Cells.Item(91, 7).Value = "Haba"
$ws.Cells.Item(91, 8).Value = "Sin especificar"
$ws.Cells.Item(91, 9).Value = "Primera"
$ws.Cells.Item(91, 10).Value = 79
$ws.Cells.Item(91, 11).Value = 6000
$ws.Cells.Item(91, 12).Value = 7000
$ws.Cells.Item(91, 13).Value = 6494
$ws.Cells.Item(91, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(91, 15).Value = "Región Metropolitana"
$ws.Cells.Item(91, 16).Value = 260
$ws.Cells.Item(91, 17).Value = 25
$ws.Cells.Item(91, 18).Value = "Hortaliza"
